$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2.3. " fragment from the annotation text in F2
$ws.Range("F2").Value = "Не-восприятие; Работа с прото-формами, не-формами, смысловые практики; Формирование объемных фигур внимания; Превращение объекта внимания в фигуру внимания;Различение действия из ментальной тишины (МТ) и пауз; Работа с локальными объемами внимания;"

# Update the selected cell to F2 as recorded in the saved workbook
$ws.Range("F2").Select()
